$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 5861
$ws.Range("J3").Value = 6275
$ws.Range("G4").Value = 1474
$ws.Range("I4").Value = 1776
$ws.Range("J4").Value = 1356
$ws.Range("J6").Value = 8026
$ws.Range("G7").Value = 24699
$ws.Range("I7").Value = 26232
$ws.Range("J7").Value = 21994

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 377
$ws.Range("J3").Value = 422
$ws.Range("J6").Value = 467
$ws.Range("J7").Value = 1379

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 245
$ws.Range("J3").Value = 342
$ws.Range("J6").Value = 351
$ws.Range("J7").Value = 1022

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 201
$ws.Range("J3").Value = 230
$ws.Range("J7").Value = 677

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J4").Value = 98
$ws.Range("J5").Value = 70
$ws.Range("J7").Value = 651
$ws.Range("J8").Value = 1379
$ws.Range("J11").Value = 346
$ws.Range("J15").Value = 245
$ws.Range("J16").Value = 86
$ws.Range("J19").Value = 653
$ws.Range("J23").Value = 207
$ws.Range("J25").Value = 110
$ws.Range("J29").Value = 1218
$ws.Range("J31").Value = 199
$ws.Range("J32").Value = 36
$ws.Range("J33").Value = 1022
$ws.Range("J36").Value = 303
$ws.Range("J37").Value = 677
$ws.Range("J42").Value = 926
$ws.Range("J46").Value = 73
$ws.Range("J47").Value = 167
$ws.Range("J48").Value = 260
$ws.Range("J49").Value = 148
$ws.Range("J50").Value = 130
$ws.Range("J51").Value = 270
$ws.Range("J52").Value = 551
$ws.Range("J54").Value = 430
$ws.Range("J55").Value = 304
$ws.Range("G63").Value = 274
$ws.Range("I63").Value = 245
$ws.Range("J67").Value = 832
$ws.Range("J68").Value = 42
$ws.Range("J76").Value = 331
$ws.Range("J77").Value = 164
$ws.Range("J79").Value = 628
$ws.Range("J85").Value = 907
$ws.Range("J86").Value = 140
$ws.Range("J88").Value = 233
$ws.Range("J91").Value = 251
$ws.Range("J94").Value = 224
$ws.Range("J96").Value = 252
$ws.Range("J97").Value = 187
$ws.Range("J98").Value = 161
$ws.Range("G101").Value = 24699
$ws.Range("I101").Value = 26232
$ws.Range("J101").Value = 21994

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J3").Value = 54
$ws.Range("J7").Value = 199

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 207
$ws.Range("J3").Value = 315
$ws.Range("J6").Value = 223
$ws.Range("J7").Value = 832

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J3").Value = 28
$ws.Range("J7").Value = 148

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J3").Value = 85
$ws.Range("J6").Value = 208
$ws.Range("J7").Value = 430

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 371
$ws.Range("J7").Value = 1218

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J6").Value = 129
$ws.Range("J7").Value = 260

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 157
$ws.Range("J3").Value = 191
$ws.Range("J6").Value = 251
$ws.Range("J7").Value = 653

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J6").Value = 184
$ws.Range("J7").Value = 331

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 190
$ws.Range("J4").Value = 42
$ws.Range("J6").Value = 478
$ws.Range("J7").Value = 926

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J6").Value = 157
$ws.Range("J7").Value = 304

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J2").Value = 22
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J3").Value = 69
$ws.Range("J7").Value = 207

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J3").Value = 69
$ws.Range("J7").Value = 252

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J6").Value = 60
$ws.Range("J7").Value = 251

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 176
$ws.Range("J6").Value = 183
$ws.Range("J7").Value = 628

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 98
$ws.Range("J3").Value = 99
$ws.Range("J7").Value = 303

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J3").Value = 195
$ws.Range("J6").Value = 211
$ws.Range("J7").Value = 651

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J2").Value = 40
$ws.Range("J6").Value = 123
$ws.Range("J7").Value = 224

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 110

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J6").Value = 79
$ws.Range("J7").Value = 167

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J4").Value = 11
$ws.Range("J7").Value = 245

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J3").Value = 24
$ws.Range("J6").Value = 99
$ws.Range("J7").Value = 161

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J6").Value = 39
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 104
$ws.Range("J7").Value = 346

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J3").Value = 21
$ws.Range("J7").Value = 187

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J4").Value = 9
$ws.Range("J6").Value = 109
$ws.Range("J7").Value = 233

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 36

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J6").Value = 33
$ws.Range("J7").Value = 70

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J4").Value = 74
$ws.Range("J6").Value = 26
$ws.Range("J7").Value = 140

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J3").Value = 72
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 270

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("J2").Value = 17
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J6").Value = 266
$ws.Range("J7").Value = 907

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J3").Value = 55
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 169
$ws.Range("J6").Value = 223
$ws.Range("J7").Value = 551

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 98

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("J6").Value = 68
$ws.Range("J7").Value = 86
